$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-16 05:48:59'
$ws.Range('G2').Value = '123 cm'
$ws.Range('H2').NumberFormat = '@'
$ws.Range('H2').Value = '99%'
$ws.Range('I2').Value = '1.9 mm'
$ws.Range('K2').Value = '0.0 MJ/m2'
$ws.Range('M2').Value = '1.8 °C 5:22 TU'
$ws.Range('N2').Value = '0.7 °C 1:29 TU'
$ws.Range('O2').Value = '1.1 °C'
$ws.Range('E3').Value = '2026-02-16 05:49:02'
$ws.Range('I3').Value = '0.5 mm'
$ws.Range('M3').Value = '-1.2 °C 5:15 TU'
$ws.Range('E4').Value = '2026-02-16 05:49:04'
$ws.Range('L4').Value = '35.6 km/h - 309º 5:29 TU'
$ws.Range('E5').Value = '2026-02-16 05:49:06'
$ws.Range('I5').Value = '1.9 mm'
$ws.Range('M5').Value = '-0.8 °C 5:27 TU'
$ws.Range('E6').Value = '2026-02-16 05:49:09'
$ws.Range('H6').NumberFormat = '@'
$ws.Range('H6').Value = '86%'
$ws.Range('N6').Value = '5.4 °C 5:13 TU'
$ws.Range('O6').Value = '6.7 °C'
$ws.Range('E7').Value = '2026-02-16 05:49:11'
$ws.Range('J7').Value = '1014.7 hPa'
$ws.Range('M7').Value = '13.8 °C 5:24 TU'
$ws.Range('E8').Value = '2026-02-16 05:49:14'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '70%'
$ws.Range('O8').Value = '9.4 °C'
$ws.Range('E9').Value = '2026-02-16 05:49:17'
$ws.Range('N9').Value = '3.8 °C 5:15 TU'
$ws.Range('O9').Value = '5.0 °C'
$ws.Range('E10').Value = '2026-02-16 05:49:19'
$ws.Range('E11').Value = '2026-02-16 05:49:22'
$ws.Range('N11').Value = '-0.2 °C 5:18 TU'
$ws.Range('O11').Value = '0.7 °C'
$ws.Range('E12').Value = '2026-02-16 05:49:24'
$ws.Range('N12').Value = '3.8 °C 5:07 TU'
$ws.Range('O12').Value = '5.4 °C'
$ws.Range('E13').Value = '2026-02-16 05:49:27'
$ws.Range('J13').Value = '1018.6 hPa'
$ws.Range('N13').Value = '-1.6 °C 5:20 TU'
$ws.Range('O13').Value = '0.9 °C'
$ws.Range('E14').Value = '2026-02-16 05:49:29'
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value = '60%'
$ws.Range('E15').Value = '2026-02-16 05:49:32'
$ws.Range('N15').Value = '3.1 °C 5:15 TU'
$ws.Range('O15').Value = '5.0 °C'
$ws.Range('E16').Value = '2026-02-16 05:49:34'
$ws.Range('M16').Value = '0.8 °C 5:01 TU'
$ws.Range('O16').Value = '-0.5 °C'
$ws.Range('E17').Value = '2026-02-16 05:49:37'
$ws.Range('E18').Value = '2026-02-16 05:49:39'
$ws.Range('H18').NumberFormat = '@'
$ws.Range('H18').Value = '99%'
$ws.Range('N18').Value = '2.8 °C 5:24 TU'
$ws.Range('O18').Value = '4.1 °C'
$ws.Range('E19').Value = '2026-02-16 05:49:42'
$ws.Range('N19').Value = '2.2 °C 5:22 TU'
$ws.Range('O19').Value = '3.2 °C'
$ws.Range('E20').Value = '2026-02-16 05:49:44'
$ws.Range('O20').Value = '-1.3 °C'
$ws.Range('E21').Value = '2026-02-16 05:49:47'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H21').Value = '82%'
$ws.Range('O21').Value = '4.6 °C'
$ws.Range('E22').Value = '2026-02-16 05:49:49'
$ws.Range('I22').Value = '0.8 mm'
$ws.Range('E23').Value = '2026-02-16 05:49:52'
$ws.Range('I23').Value = '1.0 mm'
$ws.Range('M23').Value = '-0.1 °C 5:16 TU'
$ws.Range('O23').Value = '-0.9 °C'
$ws.Range('E24').Value = '2026-02-16 05:49:54'
$ws.Range('H24').NumberFormat = '@'
$ws.Range('H24').Value = '72%'
$ws.Range('E25').Value = '2026-02-16 05:49:57'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '72%'
$ws.Range('N25').Value = '-0.6 °C 5:27 TU'
$ws.Range('O25').Value = '0.5 °C'
$ws.Range('E26').Value = '2026-02-16 05:49:59'
$ws.Range('E27').Value = '2026-02-16 05:50:02'
$ws.Range('L27').Value = '32.0 km/h - 275º 5:19 TU'
$ws.Range('E28').Value = '2026-02-16 05:50:04'
$ws.Range('N28').Value = '1.7 °C 5:02 TU'
$ws.Range('O28').Value = '3.1 °C'
$ws.Range('E29').Value = '2026-02-16 05:50:07'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '97%'
$ws.Range('N29').Value = '3.7 °C 5:29 TU'
$ws.Range('O29').Value = '4.8 °C'
$ws.Range('E30').Value = '2026-02-16 05:50:09'
$ws.Range('E31').Value = '2026-02-16 05:50:12'
$ws.Range('K31').Value = '-0.1 MJ/m2'
$ws.Range('O31').Value = '13.9 °C'
$ws.Range('E32').Value = '2026-02-16 05:50:14'
$ws.Range('H32').NumberFormat = '@'
$ws.Range('H32').Value = '82%'
$ws.Range('E33').Value = '2026-02-16 05:50:17'
$ws.Range('H33').NumberFormat = '@'
$ws.Range('H33').Value = '74%'
$ws.Range('J33').Value = '1016.2 hPa'
$ws.Range('O33').Value = '4.2 °C'
$ws.Range('E34').Value = '2026-02-16 05:50:19'
$ws.Range('E35').Value = '2026-02-16 05:50:22'
$ws.Range('E36').Value = '2026-02-16 05:50:24'
$ws.Range('J36').Value = '1014.4 hPa'
$ws.Range('L36').Value = '10.8 km/h - 87º 5:00 TU'
$ws.Range('E37').Value = '2026-02-16 05:50:27'
$ws.Range('E38').Value = '2026-02-16 05:50:30'
$ws.Range('K38').Value = '-0.1 MJ/m2'
$ws.Range('E39').Value = '2026-02-16 05:50:32'
$ws.Range('O39').Value = '0.0 °C'
$ws.Range('E40').Value = '2026-02-16 05:50:35'
$ws.Range('O40').Value = '2.8 °C'
$ws.Range('E41').Value = '2026-02-16 05:50:37'
$ws.Range('H41').NumberFormat = '@'
$ws.Range('H41').Value = '52%'
$ws.Range('J41').Value = '1015.9 hPa'
$ws.Range('E42').Value = '2026-02-16 05:50:40'
$ws.Range('H42').NumberFormat = '@'
$ws.Range('H42').Value = '95%'
$ws.Range('N42').Value = '5.4 °C 5:29 TU'
$ws.Range('E43').Value = '2026-02-16 05:50:42'
$ws.Range('O43').Value = '3.3 °C'
$ws.Range('E44').Value = '2026-02-16 05:50:44'
$ws.Range('I44').Value = '1.0 mm'
$ws.Range('E45').Value = '2026-02-16 05:50:47'
$ws.Range('I45').Value = '1.0 mm'
$ws.Range('J45').Value = '1019.5 hPa'
$ws.Range('E46').Value = '2026-02-16 05:50:49'
$ws.Range('J46').Value = '1018.5 hPa'
$ws.Range('M46').Value = '13.3 °C 5:28 TU'
